$d = $word.ActiveDocument
$p = $d.Paragraphs(27)
$full = $p.Range.WordOpenXML
Write-Output $full.Length
# try to extract the w:p...w:p element only
if ($full -match '(?s)(<w:p[ >].*?</w:p>)') {
    Write-Output "MATCHED"
    Write-Output $matches[1]
} else {
    Write-Output "no match"
}
